$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Home Decor / Home Decor
$ws.Range("C2").Value = 14248

# Row 3: Kitchen & Dining / Kitchen & Dining
$ws.Range("C3").Value = 8736

# Row 4: was Stationery & Office / Stationery & Office -> Seasonal & Holidays / Seasonal & Holidays
$ws.Range("A4").Value = "Seasonal & Holidays"
$ws.Range("B4").Value = "Seasonal & Holidays"
$ws.Range("C4").Value = 6730

# Row 5: was Seasonal & Holidays / Seasonal & Holidays -> Stationery & Office / Stationery & Office
$ws.Range("A5").Value = "Stationery & Office"
$ws.Range("B5").Value = "Stationery & Office"
$ws.Range("C5").Value = 5220

# Row 6: was Kids & Toys / Kids & Toys -> Kitchen & Dining / Home Decor
$ws.Range("A6").Value = "Kitchen & Dining"
$ws.Range("B6").Value = "Home Decor"
$ws.Range("C6").Value = 4404

# Row 7: was Textiles & Cozy Items / Textiles & Cozy Items -> Home Decor / Kitchen & Dining
$ws.Range("A7").Value = "Home Decor"
$ws.Range("B7").Value = "Kitchen & Dining"
$ws.Range("C7").Value = 4404

# Row 8: was Fashion & Accessories / Fashion & Accessories -> Seasonal & Holidays / Home Decor
$ws.Range("A8").Value = "Seasonal & Holidays"
$ws.Range("B8").Value = "Home Decor"
$ws.Range("C8").Value = 3759

# Row 9: was Party Supplies / Party Supplies -> Home Decor / Seasonal & Holidays
$ws.Range("A9").Value = "Home Decor"
$ws.Range("B9").Value = "Seasonal & Holidays"
$ws.Range("C9").Value = 3759

# Row 10: A10 stays Home Decor, B10 was Kitchen & Dining -> Stationery & Office
$ws.Range("B10").Value = "Stationery & Office"
$ws.Range("C10").Value = 2492

# Row 11: A11 was Kitchen & Dining -> Stationery & Office, B11 stays Home Decor
$ws.Range("A11").Value = "Stationery & Office"
$ws.Range("C11").Value = 2492
